$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (row 7, column B): set the Value cell to the literal
# text "true" (not the Boolean TRUE). Typing "true" directly would be
# auto-coerced to a Boolean by Excel, so build it as a text formula result
# and then paste-special as values to bake it in as a plain text cell
# (keeps the original cell style/number format intact).
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# "Date" row (row 8, column B): update the timestamp text value.
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
